$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.166777666666666
$ws.Range("H2").Value = 6.500332999999999
$ws.Range("I2").Value = 0.3769884032128669
$ws.Range("J2").Value = 0.376988403212867
$ws.Range("M2").Value = 39.62362533333334
$ws.Range("N2").Value = 118.870876
$ws.Range("O2").Value = 0.862331223864382
$ws.Range("P2").Value = 0.8623312238643819
$ws.Range("Q2").Value = 85.85558644463423
$ws.Range("R2").Value = 772.7002780017081
$ws.Range("S2").Value = 0.3250888711252307
$ws.Range("T2").Value = 0.3250888711252307
$ws.Range("G3").Value = 2.166777666666666
$ws.Range("H3").Value = 6.500332999999999
$ws.Range("I3").Value = 0.3769884032128669
$ws.Range("J3").Value = 0.376988403212867
$ws.Range("O3").Value = 0.04891146266025228
$ws.Range("P3").Value = 0.04891146266025227
$ws.Range("Q3").Value = 4.869732411801444
$ws.Range("R3").Value = 43.827591706213
$ws.Range("S3").Value = 0.01843905420709427
$ws.Range("T3").Value = 0.01843905420709427
$ws.Range("G4").Value = 2.166777666666666
$ws.Range("H4").Value = 6.500332999999999
$ws.Range("I4").Value = 0.3769884032128669
$ws.Range("J4").Value = 0.376988403212867
$ws.Range("M4").Value = 0.2964306666666667
$ws.Range("N4").Value = 0.889292
$ws.Range("O4").Value = 0.006451237548992269
$ws.Range("P4").Value = 0.006451237548992268
$ws.Range("Q4").Value = 0.6422993482484444
$ws.Range("R4").Value = 5.780694134236
$ws.Range("S4").Value = 0.002432041742341485
$ws.Range("T4").Value = 0.002432041742341485
$ws.Range("G5").Value = 2.166777666666666
$ws.Range("H5").Value = 6.500332999999999
$ws.Range("I5").Value = 0.3769884032128669
$ws.Range("J5").Value = 0.376988403212867
$ws.Range("M5").Value = 2.447182
$ws.Range("N5").Value = 7.341546
$ws.Range("O5").Value = 0.05325816179933475
$ws.Range("P5").Value = 0.05325816179933474
$ws.Range("Q5").Value = 5.302499303868666
$ws.Range("R5").Value = 47.72249373481799
$ws.Range("S5").Value = 0.02007770937478371
$ws.Range("T5").Value = 0.02007770937478371
$ws.Range("G6").Value = 2.166777666666666
$ws.Range("H6").Value = 6.500332999999999
$ws.Range("I6").Value = 0.3769884032128669
$ws.Range("J6").Value = 0.376988403212867
$ws.Range("M6").Value = 0.4200656666666667
$ws.Range("N6").Value = 1.260197
$ws.Range("O6").Value = 0.009141913123616776
$ws.Range("P6").Value = 0.009141913123616775
$ws.Range("Q6").Value = 0.9101889050667777
$ws.Range("R6").Value = 8.191700145600999
$ws.Range("S6").Value = 0.003446395230783041
$ws.Range("T6").Value = 0.003446395230783041
$ws.Range("G7").Value = 2.166777666666666
$ws.Range("H7").Value = 6.500332999999999
$ws.Range("I7").Value = 0.3769884032128669
$ws.Range("J7").Value = 0.376988403212867
$ws.Range("M7").Value = 0.9146693333333333
$ws.Range("N7").Value = 2.744008
$ws.Range("O7").Value = 0.01990600100342202
$ws.Range("P7").Value = 0.01990600100342202
$ws.Range("Q7").Value = 1.981885083851555
$ws.Range("R7").Value = 17.836965754664
$ws.Range("S7").Value = 0.007504331532633795
$ws.Range("T7").Value = 0.007504331532633795
$ws.Range("I8").Value = 0.3757968909097267
$ws.Range("J8").Value = 0.3757968909097268
$ws.Range("M8").Value = 39.62362533333334
$ws.Range("N8").Value = 118.870876
$ws.Range("O8").Value = 0.862331223864382
$ws.Range("P8").Value = 0.8623312238643819
$ws.Range("Q8").Value = 85.58423065047646
$ws.Range("R8").Value = 770.2580758542881
$ws.Range("S8").Value = 0.3240613928626143
$ws.Range("T8").Value = 0.3240613928626143
$ws.Range("I9").Value = 0.3757968909097267
$ws.Range("J9").Value = 0.3757968909097268
$ws.Range("O9").Value = 0.04891146266025228
$ws.Range("P9").Value = 0.04891146266025227
$ws.Range("S9").Value = 0.01838077559757
$ws.Range("T9").Value = 0.01838077559757
$ws.Range("I10").Value = 0.3757968909097267
$ws.Range("J10").Value = 0.3757968909097268
$ws.Range("M10").Value = 0.2964306666666667
$ws.Range("N10").Value = 0.889292
$ws.Range("O10").Value = 0.006451237548992269
$ws.Range("P10").Value = 0.006451237548992268
$ws.Range("Q10").Value = 0.6402692922328889
$ws.Range("R10").Value = 5.762423630096
$ws.Range("S10").Value = 0.002424355013431381
$ws.Range("T10").Value = 0.002424355013431381
$ws.Range("I11").Value = 0.3757968909097267
$ws.Range("J11").Value = 0.3757968909097268
$ws.Range("M11").Value = 2.447182
$ws.Range("N11").Value = 7.341546
$ws.Range("O11").Value = 0.05325816179933475
$ws.Range("P11").Value = 0.05325816179933474
$ws.Range("Q11").Value = 5.285740185805333
$ws.Range("R11").Value = 47.571661672248
$ws.Range("S11").Value = 0.02001425161975718
$ws.Range("T11").Value = 0.02001425161975718
$ws.Range("I12").Value = 0.3757968909097267
$ws.Range("J12").Value = 0.3757968909097268
$ws.Range("M12").Value = 0.4200656666666667
$ws.Range("N12").Value = 1.260197
$ws.Range("O12").Value = 0.009141913123616776
$ws.Range("P12").Value = 0.009141913123616775
$ws.Range("Q12").Value = 0.9073121553595556
$ws.Range("R12").Value = 8.165809398236
$ws.Range("S12").Value = 0.003435502528822013
$ws.Range("T12").Value = 0.003435502528822013
$ws.Range("I13").Value = 0.3757968909097267
$ws.Range("J13").Value = 0.3757968909097268
$ws.Range("M13").Value = 0.9146693333333333
$ws.Range("N13").Value = 2.744008
$ws.Range("O13").Value = 0.01990600100342202
$ws.Range("P13").Value = 0.01990600100342202
$ws.Range("Q13").Value = 1.975621123367111
$ws.Range("R13").Value = 17.780590110304
$ws.Range("S13").Value = 0.007480613287531896
$ws.Range("T13").Value = 0.007480613287531896
$ws.Range("G14").Value = 1.420890666666667
$ws.Range("H14").Value = 4.262672
$ws.Range("I14").Value = 0.2472147058774063
$ws.Range("J14").Value = 0.2472147058774063
$ws.Range("M14").Value = 39.62362533333334
$ws.Range("N14").Value = 118.870876
$ws.Range("O14").Value = 0.862331223864382
$ws.Range("P14").Value = 0.8623312238643819
$ws.Range("Q14").Value = 56.30083941563023
$ws.Range("R14").Value = 506.7075547406721
$ws.Range("S14").Value = 0.213180959876537
$ws.Range("T14").Value = 0.213180959876537
$ws.Range("G15").Value = 1.420890666666667
$ws.Range("H15").Value = 4.262672
$ws.Range("I15").Value = 0.2472147058774063
$ws.Range("J15").Value = 0.2472147058774063
$ws.Range("O15").Value = 0.04891146266025228
$ws.Range("P15").Value = 0.04891146266025227
$ws.Range("Q15").Value = 3.193385938732445
$ws.Range("R15").Value = 28.74047344859201
$ws.Range("S15").Value = 0.01209163285558801
$ws.Range("T15").Value = 0.01209163285558801
$ws.Range("G16").Value = 1.420890666666667
$ws.Range("H16").Value = 4.262672
$ws.Range("I16").Value = 0.2472147058774063
$ws.Range("J16").Value = 0.2472147058774063
$ws.Range("M16").Value = 0.2964306666666667
$ws.Range("N16").Value = 0.889292
$ws.Range("O16").Value = 0.006451237548992269
$ws.Range("P16").Value = 0.006451237548992268
$ws.Range("Q16").Value = 0.4211955675804445
$ws.Range("R16").Value = 3.790760108224
$ws.Range("S16").Value = 0.001594840793219403
$ws.Range("T16").Value = 0.001594840793219403
$ws.Range("G17").Value = 1.420890666666667
$ws.Range("H17").Value = 4.262672
$ws.Range("I17").Value = 0.2472147058774063
$ws.Range("J17").Value = 0.2472147058774063
$ws.Range("M17").Value = 2.447182
$ws.Range("N17").Value = 7.341546
$ws.Range("O17").Value = 0.05325816179933475
$ws.Range("P17").Value = 0.05325816179933474
$ws.Range("Q17").Value = 3.477178063434667
$ws.Range("R17").Value = 31.294602570912
$ws.Range("S17").Value = 0.01316620080479386
$ws.Range("T17").Value = 0.01316620080479386
$ws.Range("G18").Value = 1.420890666666667
$ws.Range("H18").Value = 4.262672
$ws.Range("I18").Value = 0.2472147058774063
$ws.Range("J18").Value = 0.2472147058774063
$ws.Range("M18").Value = 0.4200656666666667
$ws.Range("N18").Value = 1.260197
$ws.Range("O18").Value = 0.009141913123616776
$ws.Range("P18").Value = 0.009141913123616775
$ws.Range("Q18").Value = 0.5968673851537778
$ws.Range("R18").Value = 5.371806466384
$ws.Range("S18").Value = 0.002260015364011722
$ws.Range("T18").Value = 0.002260015364011722
$ws.Range("G19").Value = 1.420890666666667
$ws.Range("H19").Value = 4.262672
$ws.Range("I19").Value = 0.2472147058774063
$ws.Range("J19").Value = 0.2472147058774063
$ws.Range("M19").Value = 0.9146693333333333
$ws.Range("N19").Value = 2.744008
$ws.Range("O19").Value = 0.01990600100342202
$ws.Range("P19").Value = 0.01990600100342202
$ws.Range("Q19").Value = 1.299645118819556
$ws.Range("R19").Value = 11.696806069376
$ws.Range("S19").Value = 0.00492105618325633
$ws.Range("T19").Value = 0.004921056183256329
